$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily status row (row 50) below the existing data (row 49).
$row = 50

$ws.Cells.Item($row, 1).Value = 46007
$ws.Cells.Item($row, 2).Value = 5610
$ws.Cells.Item($row, 3).Value = 4328
$ws.Cells.Item($row, 4).Value = 4016
$ws.Cells.Item($row, 5).Value = 226
$ws.Cells.Item($row, 6).Value = 46
$ws.Cells.Item($row, 7).Value = 35
$ws.Cells.Item($row, 8).Value = 5
$ws.Cells.Item($row, 9).Value = 0

# Match the date formatting used by the "Date" column in the row above it.
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

# Reflect the new active selection on the sheet (row 50 across columns A:I).
$null = $ws.Range("A50:I50").Select()

Write-Host "Added row 50 to Sheet1"
